# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled
# update). For every data row, Price (D) and Volume(1h) (E) get the
# latest scraped figures; rows 47/48 additionally swap which coin
# (EnergySwap / Frax) occupies which rank, so Coin (B) and Link (C)
# are rewritten too.
#
# Column D holds plain text in the source sheet (e.g. "26.290.80",
# "0.5262") even when a value happens to parse as a number. Writing a
# numeric-looking string straight into `.Value` lets Excel's COM layer
# silently coerce it to a float (losing trailing zeros / exact
# formatting), so Set-TextValue forces those via the classic
# leading-apostrophe "text" prefix and then resets the cell style to
# Normal so no stray NumberFormat/quote-prefix formatting is left
# behind. Column E values never look numeric (they keep padding
# spaces, e.g. "  +1.08%  "), so they're assigned directly.

function Set-TextValue($ws, $cellRef, $val) {
    $ws.Range($cellRef).Value = "'" + $val
    $ws.Range($cellRef).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range('D2').Value = '26.303.20'
$ws.Range('E2').Value = '  +1.08%  '
$ws.Range('D3').Value = '1.678.06'
$ws.Range('E3').Value = '  +0.78%  '
$ws.Range('E4').Value = '  +0.33%  '
$ws.Range('E5').Value = '  +0.72%  '
Set-TextValue $ws 'D6' '0.5259'
$ws.Range('E6').Value = '  +3.35%  '
Set-TextValue $ws 'D7' '1.009'
Set-TextValue $ws 'D8' '0.2685'
$ws.Range('E8').Value = '  +2.17%  '
Set-TextValue $ws 'D9' '0.06461'
$ws.Range('E9').Value = '  +1.08%  '
Set-TextValue $ws 'D10' '21.84'
$ws.Range('E10').Value = '  +0.70%  '
Set-TextValue $ws 'D11' '0.07511'
$ws.Range('E11').Value = '  +1.26%  '
$ws.Range('D12').Value = '1.690.40'
$ws.Range('E12').Value = '  +1.29%  '
Set-TextValue $ws 'D13' '4.509'
$ws.Range('E13').Value = '  +0.23%  '
Set-TextValue $ws 'D14' '0.5772'
$ws.Range('E14').Value = '  -0.57%  '
Set-TextValue $ws 'D15' '0.000008496'
$ws.Range('E15').Value = '  -0.10%  '
Set-TextValue $ws 'D16' '64.65'
$ws.Range('E16').Value = '  +0.75%  '
$ws.Range('D17').Value = '26.331.07'
$ws.Range('E17').Value = '  +0.91%  '
Set-TextValue $ws 'D18' '4.906'
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('E19').Value = '  +0.27%  '
Set-TextValue $ws 'D20' '10.86'
$ws.Range('E20').Value = '  +1.61%  '
Set-TextValue $ws 'D21' '189.50'
$ws.Range('E21').Value = '  +0.48%  '
Set-TextValue $ws 'D22' '6.175'
$ws.Range('E22').Value = '  -0.37%  '
$ws.Range('E23').Value = '  +0.25%  '
Set-TextValue $ws 'D24' '144.91'
$ws.Range('E24').Value = '  -0.55%  '
Set-TextValue $ws 'D25' '7.768'
$ws.Range('E25').Value = '  +2.20%  '
Set-TextValue $ws 'D26' '0.1258'
$ws.Range('E26').Value = '  +5.94%  '
Set-TextValue $ws 'D27' '15.75'
$ws.Range('E27').Value = '  +0.96%  '
$ws.Range('E28').Value = '  -3.18%  '
Set-TextValue $ws 'D29' '1.364'
$ws.Range('E29').Value = '  +4.52%  '
Set-TextValue $ws 'D30' '1.323'
$ws.Range('E30').Value = '  +0.81%  '
Set-TextValue $ws 'D31' '3.579'
$ws.Range('E31').Value = '  +1.63%  '
Set-TextValue $ws 'D32' '3.587'
$ws.Range('E32').Value = '  +2.52%  '
Set-TextValue $ws 'D33' '1.655'
$ws.Range('E33').Value = '  +1.81%  '
Set-TextValue $ws 'D34' '1.027'
$ws.Range('E34').Value = '  +0.92%  '
Set-TextValue $ws 'D35' '0.6191'
$ws.Range('E35').Value = '  +2.31%  '
Set-TextValue $ws 'D36' '2.406'
$ws.Range('E36').Value = '  +1.64%  '
Set-TextValue $ws 'D37' '2.742'
$ws.Range('E37').Value = '  +2.19%  '
Set-TextValue $ws 'D38' '6.278'
$ws.Range('E38').Value = '  +1.28%  '
$ws.Range('D39').Value = '1.115.65'
$ws.Range('E39').Value = '  +3.89%  '
$ws.Range('E40').Value = '  +0.75%  '
Set-TextValue $ws 'D41' '0.8710'
$ws.Range('E41').Value = '  +1.38%  '
Set-TextValue $ws 'D42' '1.016'
$ws.Range('E42').Value = '  +0.70%  '
Set-TextValue $ws 'D43' '100.43'
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('D44').Value = '1.828.53'
$ws.Range('E44').Value = '  +0.89%  '
$ws.Range('E45').Value = '  -5.86%  '
Set-TextValue $ws 'D46' '56.86'
$ws.Range('E46').Value = '  +1.19%  '
$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws 'D47' '1.006'
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws 'D48' '8.131'
$ws.Range('E48').Value = '  +1.75%  '
Set-TextValue $ws 'D49' '0.05265'
$ws.Range('E49').Value = '  +1.08%  '
Set-TextValue $ws 'D50' '0.4298'
$ws.Range('E50').Value = '  +0.12%  '
Set-TextValue $ws 'D51' '6.047'
$ws.Range('E51').Value = '  +1.65%  '
